$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, shifting existing row 151 (and below) down to 152.
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with its data.
$ws.Cells.Item(151, 1).Value = 3
$ws.Cells.Item(151, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44510
$ws.Cells.Item(151, 5).Value = 5
$ws.Cells.Item(151, 6).Value = 100112012
$ws.Cells.Item(151, 7).Value = "Espinaca"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 130
$ws.Cells.Item(151, 11).Value = 2000
$ws.Cells.Item(151, 12).Value = 2500
$ws.Cells.Item(151, 13).Value = 2231
$ws.Cells.Item(151, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(151, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(151, 16).Value = 744
$ws.Cells.Item(151, 17).Value = 3
$ws.Cells.Item(151, 18).Value = "Hortaliza"
